# Atualizacao DDD Turma 2
# Updates the "Plano de Aula DDD" sheet for the new class (Turma 2):
#  - Window is no longer minimized
#  - C2: clear the leftover "Turma 1" notes cell (content removed, wrap-text style kept)
#  - C3: clear the leftover "Turma 1" notes cell entirely (cell removed)
#  - Selection moves to A13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plano de Aula DDD")

# restore the window to a normal (non-minimized) state
$win = $wb.Windows.Item(1)
$win.WindowState = -4143

$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

$ws.Range("A13").Select()
